{"js": "// \"Simple of ETL on [Application].[People]:\" -> \"Sample of ETL on [Application].[People]:\"\n//   - the fixed word is typed back in as three separate runs: \"S\" / \"a\" / \"mple of ETL on \"\n//   - the \"[People]\" run and the \":\" run (identical formatting) become a single \"[People]:\" run\nconst body = context.document.body;\n\n// ---- Edit 1: \"Simple\" -> \"Sample\", split across three runs: \"S\", \"a\", \"mple of ETL on \" ----\nconst fullResults = body.search(\"Simple of ETL on \", { matchCase: true });\nfullResults.load(\"items\");\nawait context.sync();\nconst fullRange = fullResults.items[0];\n\n// Collapsed anchor at the very start of \"Simple of ETL on \", then clear the old text.\nconst startRange = fullRange.getRange(\"Start\");\nfullRange.insertText(\"\", \"Replace\");\nawait context.sync();\n\n// Type the pieces back in right-to-left, re-locating each insertion point with a fresh\n// search. Typing left-to-right (or reusing a stale collapsed Range) lets the engine\n// coalesce the newly-typed text with its same-formatted neighbour into one run; inserting\n// \"Before\" a freshly re-found point keeps each piece in its own run.\nstartRange.insertText(\"mple of ETL on \", \"Before\");\nawait context.sync();\n\nconst mResults = body.search(\"mple of ETL on \", { matchCase: true });\nmResults.load(\"items\");\nawait context.sync();\nmResults.items[0].getRange(\"Start\").insertText(\"a\", \"Before\");\nawait context.sync();\n\nconst aResults = body.search(\"ample of ETL on \", { matchCase: true });\naResults.load(\"items\");\nawait context.sync();\naResults.items[0].getRange(\"Start\").insertText(\"S\", \"Before\");\nawait context.sync();\n\n// ---- Edit 2: merge the \"[People]\" run and the \":\" run into a single \"[People]:\" run ----\nconst peopleResults = body.search(\"[People]:\", { matchCase: true });\npeopleResults.load(\"items\");\nawait context.sync();\npeopleResults.items[0].insertText(\"[People]:\", \"Replace\");\nawait context.sync();\n", "ps1": "# \"Simple of ETL on [Application].[People]:\" -> \"Sample of ETL on [Application].[People]:\"\n#   - the fixed word is typed back in as three separate runs: \"S\" / \"a\" / \"mple of ETL on \"\n#   - the \"[People]\" run and the \":\" run (identical formatting) become a single \"[People]:\" run\n$d = $word.ActiveDocument\n\n# ---- Edit 1: \"Simple\" -> \"Sample\", split across three runs: \"S\", \"a\", \"mple of ETL on \" ----\n$find = $d.Content.Find\n$find.Text = \"Simple of ETL on \"\n$find.Execute() | Out-Null\n$fullRange = $find.Parent\n\n# Collapsed anchor at the very start of \"Simple of ETL on \", then clear the old text.\n$startRange = $d.Range($fullRange.Start, $fullRange.Start)\n$fullRange.Text = \"\"\n\n# Type the pieces back in right-to-left, re-locating each insertion point with a fresh\n# Find each time. Typing left-to-right (or reusing a stale collapsed Range) lets the\n# engine coalesce the newly-typed text with its same-formatted neighbour into one run;\n# inserting before a freshly re-found point keeps each piece in its own run.\n$startRange.InsertBefore(\"mple of ETL on \")\n\n$find2 = $d.Content.Find\n$find2.Text = \"mple of ETL on \"\n$find2.Execute() | Out-Null\n$mRange = $find2.Parent\n$beforeM = $d.Range($mRange.Start, $mRange.Start)\n$beforeM.InsertBefore(\"a\")\n\n$find3 = $d.Content.Find\n$find3.Text = \"ample of ETL on \"\n$find3.Execute() | Out-Null\n$aRange = $find3.Parent\n$beforeA = $d.Range($aRange.Start, $aRange.Start)\n$beforeA.InsertBefore(\"S\")\n\n# ---- Edit 2: merge the \"[People]\" run and the \":\" run into a single \"[People]:\" run ----\n$find4 = $d.Content.Find\n$find4.Text = \"[People]:\"\n$find4.Execute() | Out-Null\n$peopleRange = $find4.Parent\n$peopleRange.Text = \"[People]:\"\n"}
